$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.363.18"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'226.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'32.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "'0.295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'0.0690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.784.07"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'11.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "'0.633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "34.374.74"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "'68.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "0.0₃0796"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'243.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'11.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "'165.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "'7.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'3.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.22%  "
$ws.Range("D31").Value = "'0.0524"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "'3.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").Value = "1.399.61"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").Value = "'0.673"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'1.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'84.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("D44").Value = "'13.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "'1.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("D47").Value = "'6.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "1.948.21"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "'104.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  -1.67%  "
